$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "311.83"
Set-TextValue "E2" "-0.51%"
Set-TextValue "D3" "48.28"
Set-TextValue "E3" "8.03%"
Set-TextValue "D4" "5.250"
Set-TextValue "E4" "2.32%"
Set-TextValue "D5" "0.07832"
Set-TextValue "E5" "-2.63%"
Set-TextValue "D6" "4.565"
Set-TextValue "E6" "1.04%"
Set-TextValue "D7" "1.327"
Set-TextValue "E7" "22.46%"
Set-TextValue "D8" "1.591"
Set-TextValue "E8" "-5.96%"
Set-TextValue "D9" "0.1253"
Set-TextValue "E9" "-3.30%"
Set-TextValue "D10" "0.1961"
Set-TextValue "E10" "2.22%"
Set-TextValue "D11" "0.09369"
Set-TextValue "E11" "-0.24%"
Set-TextValue "D12" "0.04550"
Set-TextValue "E12" "7.49%"
Set-TextValue "D13" "0.1045"
Set-TextValue "E13" "0.35%"
Set-TextValue "D14" "0.001298"
Set-TextValue "E14" "-1.03%"
Set-TextValue "D15" "0.04203"
Set-TextValue "E15" "-0.03%"
Set-TextValue "D16" "0.005838"
Set-TextValue "E16" "-1.42%"
Set-TextValue "D17" "3.336"
Set-TextValue "E17" "-1.65%"
Set-TextValue "D18" "2.428"
Set-TextValue "E18" "1.04%"
Set-TextValue "E19" "2.07%"
Set-TextValue "D20" "8.120"
Set-TextValue "E20" "0.86%"
Set-TextValue "D21" "0.1366"
Set-TextValue "E21" "-0.29%"
Set-TextValue "D22" "0.3069"
Set-TextValue "E22" "-2.17%"
Set-TextValue "D23" "0.001295"
Set-TextValue "E23" "1.84%"
Set-TextValue "D24" "0.004199"
Set-TextValue "E24" "-8.03%"
Set-TextValue "D25" "0.0001357"
Set-TextValue "E25" "1.17%"
Set-TextValue "D26" "0.0003557"
Set-TextValue "E26" "-95.20%"
Set-TextValue "D38" "0.02598"
Set-TextValue "E38" "-4.30%"
Set-TextValue "D39" "0.05806"
Set-TextValue "E39" "7.07%"
Set-TextValue "D40" "0.01082"
Set-TextValue "E40" "91.94%"
Set-TextValue "E41" "3.28%"
Set-TextValue "E42" "1.84%"
Set-TextValue "D43" "0.008410"
Set-TextValue "E43" "14.69%"
Set-TextValue "D44" "0.008609"
Set-TextValue "E44" "8.31%"
Set-TextValue "D45" "0.3131"
Set-TextValue "E45" "0.02%"
Set-TextValue "D46" "0.00006941"
Set-TextValue "E46" "2.18%"
Set-TextValue "D47" "0.00000000754"
Set-TextValue "E47" "1.14%"
Set-TextValue "E48" "-20.93%"
Set-TextValue "D49" "0.004020"
Set-TextValue "E49" "1.13%"
Set-TextValue "D50" "0.00002110"
Set-TextValue "E50" "1.14%"
Set-TextValue "D51" "0.0002010"
Set-TextValue "E51" "1.14%"
